# Qatar Stars League base update (01-06-2024 01:16)
# The source data for several fixtures was re-ordered; for each pair of rows
# listed below, the data in columns B:AD (everything except the running
# index in column A) needs to be swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(15, 16),
    @(68, 69),
    @(96, 97),
    @(100, 101),
    @(104, 105),
    @(108, 109),
    @(118, 119),
    @(122, 123),
    @(124, 125),
    @(128, 129)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
